$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the date format used by the existing "Fecha" column (D) so that
# the new row we insert keeps the exact same style as the rest of the sheet.
$dateFormat = $ws.Range("D14").NumberFormat()

# --- Step 1: move the current (old) row 14 record down to a new row 15 ---
# Row 15 becomes an exact copy of the original row 14 values (same date,
# same volume/price figures) since that week's record is preserved as
# historical data once the new week's record is added to row 14.
$ws.Range("A15").Value = 8
$ws.Range("B15").Value = "Terminal La Palmera de La Serena"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 44466
$ws.Range("D15").NumberFormat = $dateFormat
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100101
$ws.Range("H15").Value = "Berries"
$ws.Range("I15").Value = 100101001
$ws.Range("J15").Value = "Arándano (blue)"
$ws.Range("K15").Value = "Sin especificar"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 160
$ws.Range("N15").Value = 13500
$ws.Range("O15").Value = 14000
$ws.Range("P15").Value = 13750
$ws.Range("Q15").Value = "$/bandeja 2 kilos"
$ws.Range("R15").Value = "Provincia de Limarí"
$ws.Range("S15").Value = 6875
$ws.Range("T15").Value = 2

# --- Step 2: overwrite row 14 with the new week's record ---
$ws.Range("D14").Value = 44491
$ws.Range("D14").NumberFormat = $dateFormat
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 11500
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 11750
$ws.Range("S14").Value = 5875
$ws.Range("T14").Value = 2
